# Updates license text on all "footer" sheets from "CC BY" to "CC BY SA".
#
# For each slide, find the footer shape whose name starts with
# "RStudio® is a trademark" (shapes 695 and 924 in the original deck) and:
#   1. Update the run containing "CC BY " -> "CC BY SA" (drop trailing space,
#      keep the hyperlink run/formatting intact).
#   2. Update the following run "RStudio \u2022  " by prepending the space(s)
#      that were dropped from the previous run, so the rendered text is
#      unchanged apart from "BY" -> "BY SA".
#   3. Update the shape's Name (cNvPr/@name) to match, since PowerPoint keeps
#      the shape Name in sync with the autogenerated text summary.
#
# Note: the two occurrences in the source deck differ by one space in the
# resulting gap before "RStudio" (a quirk of the original authored edit),
# so the padding is looked up per shape Id to reproduce it exactly.

$p = $ppt.ActivePresentation

# NOTE: this host's TextRange/Characters getter mangles the bullet
# character (U+2022 "\u2022") into a plain "*" when its .Text is read back
# as a string. Round-tripping a run's .Text through a read + rewrite would
# therefore silently corrupt it, so the run that contains the bullet
# ("RStudio \u2022  ") is replaced with a freshly-built literal string
# (using [char]0x2022) instead of reusing the value read from the shape.
$bullet = [char]0x2022

# Shape Id -> literal whitespace to place at the start of the run that used
# to read "RStudio \u2022  " (immediately following the "CC BY " run).
$paddingById = @{ 695 = "  "; 924 = " " }
$defaultPadding = " "

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if (-not $shape.HasTextFrame) {
            continue
        }
        if ($shape.Name -notlike "RStudio*is a trademark*") {
            continue
        }

        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text

        $ccIndex = $fullText.IndexOf("CC BY ")
        if ($ccIndex -lt 0) {
            continue
        }

        # The text box uses spAutoFit; editing its runs makes the host
        # recompute an "ideal" height from current font metrics, which can
        # drift slightly from the author's original size. Preserve the
        # original height explicitly so only the text/name actually change
        # (Left/Top/Width are left untouched to avoid needless point/EMU
        # round-trip rounding).
        $origHeight = $shape.Height

        if ($paddingById.ContainsKey($shape.Id)) {
            $padding = $paddingById[$shape.Id]
        } else {
            $padding = $defaultPadding
        }

        # 1-based start position of the "CC BY " run.
        $ccStart = $ccIndex + 1
        $ccOldLen = 6
        $ccRun = $tr.Characters($ccStart, $ccOldLen)
        $ccNewText = "CC BY SA"
        $ccRun.Text = $ccNewText

        # The run right after "CC BY " (originally "RStudio \u2022  ").
        # Recompute the start using the *new* length of the run we just
        # replaced, since the text shifted after the assignment above.
        # Rebuild the literal text instead of reading-and-reusing .Text,
        # to avoid the bullet-mangling getter bug noted above.
        $afterStart = $ccStart + $ccNewText.Length
        $afterLen = 11
        $afterRun = $tr.Characters($afterStart, $afterLen)
        $afterRun.Text = $padding + "RStudio " + $bullet + "  "

        $shape.Name = $shape.Name.Replace("CC BY RStudio", "CC BY SA" + $padding + "RStudio")

        # Restore the original height (see note above).
        $shape.Height = $origHeight
    }
}
